$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (top-right block). Written in this specific order so the
# new shared strings are appended to sharedStrings.xml in the same order as
# the target workbook (entrega, comunicación, revisión, revisiones).
$ws.Range("E1").Value = "Fecha de entrega:"
$ws.Range("A13").Value = "Encargado de comunicación"
$ws.Range("D1").Value = "Fecha de revisión:"
$ws.Range("A12").Value = "Encargado de revisiones"

# New column widths for the added columns B:E (target widths, in Excel's
# stored character units, are 13.140625 / 11.140625 / 17 / 16.7109375; the
# engine snaps ColumnWidth to its internal pixel grid, so the values below
# are the inputs that land closest - exactly, where achievable - on those
# stored widths once read back from the saved OOXML).
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 10.333333333333334
$ws.Columns.Item(4).ColumnWidth = 16.166666666666668
$ws.Columns.Item(5).ColumnWidth = 15.833333333333334

# Update the active selection to A5.
$ws.Range("A5").Select()
